# Add a new "Json.Net" license row (row 16) to the license overview sheet,
# mirroring the formatting of the row above it (row 15) and adding a
# hyperlink for the new "Link" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 15) into the
# new row (row 16) by copying the whole row range across.
$ws.Range("A15:D15").Copy($ws.Range("A16:D16"))

# Match the auto-computed row height used by the other wrapped-text rows.
$ws.Rows.Item(16).RowHeight = 28.8

# Fill in the new library's data.
$ws.Range("A16").Value = "Json.Net"
$ws.Range("B16").Value = "MIT"
$ws.Range("C16").Value = "http://www.newtonsoft.com/json"
$ws.Range("D16").Value = " - Library for serialization/deserialization into/from json format`n - Included as Dll, code is available on the project's homepage"

# Turn the link cell into a real hyperlink, like the other rows.
$ws.Hyperlinks.Add($ws.Range("C16"), "http://www.newtonsoft.com/json")

# Adding the hyperlink re-applies Excel's built-in "Link" cell style, which
# would overwrite the table formatting copied above, so restore just the
# formats (matching the other hyperlink cells, e.g. C15) afterwards.
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# Reflect the cursor position that was left selected when the workbook was
# last saved.
$null = $ws.Range("D13").Select()
